# Update column G ("K") values for rows 2-14 on the active sheet.
# This reflects regenerating save_data to use K (strikeouts) instead of
# Strike# for the pitcher, with the new, recalculated values below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 4
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
